# Add "Richard Turgeon" to the Argonne team list (slide 5, "Team" slide),
# right before "Justin Wozniak", and drop the stray empty bullet-less
# paragraph that used to sit right after that same group (after "Harry Yoo").

$p = $ppt.ActivePresentation

# Locate the "Team" slide / content placeholder robustly by scanning for the
# shape whose text contains "Justin Wozniak" rather than hard-coding a slide
# number.
$targetSlide = $null
$targetShape = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            $txt = $shape.TextFrame.TextRange.Text
            if ($txt -ne $null -and $txt.Contains("Justin Wozniak")) {
                $targetSlide = $slide
                $targetShape = $shape
                break
            }
        }
    }
    if ($targetShape -ne $null) { break }
}

$tr = $targetShape.TextFrame.TextRange

# --- 1. Insert a new "Richard Turgeon" paragraph right before "Justin Wozniak" ---

$justin = $tr.Find("Justin Wozniak")
$justinStart = $justin.Start

$paraCount = $tr.Paragraphs().Count
$justinIdx = -1
for ($i = 1; $i -le $paraCount; $i++) {
    $para = $tr.Paragraphs($i, 1)
    if ($para.Start -eq $justinStart) {
        $justinIdx = $i
        break
    }
}

$justinPara = $tr.Paragraphs($justinIdx, 1)
$justinPara.InsertBefore("Richard Turgeon`r") | Out-Null

# --- 2. Remove the empty, bullet-less paragraph right after the "Harry Yoo" /
#        "..." pair (it used to separate the Argonne block from "Fredrick"). ---

$harry = $tr.Find("Harry")
$harryStart = $harry.Start

$paraCount = $tr.Paragraphs().Count
$harryIdx = -1
for ($i = 1; $i -le $paraCount; $i++) {
    $para = $tr.Paragraphs($i, 1)
    if ($para.Start -eq $harryStart) {
        $harryIdx = $i
        break
    }
}

# Harry Yoo paragraph -> "..." paragraph -> blank paragraph
$blankPara = $tr.Paragraphs($harryIdx + 2, 1)
if ($blankPara.Text.Trim() -eq "") {
    $blankPara.Delete() | Out-Null
}
